# Add U21, a 1:4 clock fanout buffer (NB3N551, SOIC8), to the BOM.
# It is inserted right before the existing row for X1 (the USB connector),
# which was row 98 and becomes row 99.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Make room for the new part by inserting a blank row at 98; this shifts
# every row from 98 onward down by one and auto-updates dependent formulas
# (the K2 grand-total SUM range and the dimension/used-range).
$ws.Rows.Item(98).Insert()

# Fill in the new BOM row for U21.
$ws.Range("A98").Value = "U21"
$ws.Range("B98").Value = "NB3N551"
$ws.Range("C98").Value = "NB3N551"
$ws.Range("D98").Value = "SOIC8"
$ws.Range("E98").Value = "DK"
$ws.Range("F98").Value = "NB3N551DGOS-ND"
$ws.Range("G98").Value = "NB3N551DG"
$ws.Range("H98").Value = 1
$ws.Range("I98").Value = 1.86
$ws.Range("J98").Formula = "=H98*I98"

# Keep the trailing part-price formula (now on row 102, was row 101) in the
# same shared-formula style as its neighbours.
$ws.Range("J102").Formula = "=H102*I102"

# Restore the active cell / selection reflected in the saved view state.
$ws.Range("I97").Select()
